$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# --- Update path values (column B) for the existing Folder path rows ---
# Row 13: jobFolderPath stays the same name, but now points to the Archive folder
$ws.Cells.Item(13, 2).Value = "Data\Archive\Jobs"
$ws.Cells.Item(13, 1).Value = "jobFolderPath"

# Row 14: cvFolderPath stays the same name, but now points to the Archive folder
$ws.Cells.Item(14, 2).Value = "Data\Archive\CVs"
$ws.Cells.Item(14, 1).Value = "cvFolderPath"

# Row 15: repurposed from jobIndexPath -> employeesListsPath
$ws.Cells.Item(15, 2).Value = "Data\Employee"

# Row 19 value (email body for the new employees list notification)
$ws.Cells.Item(19, 2).Value = "Hi! This is the list of potential future employees for the job with the id: "

# Row 15 name
$ws.Cells.Item(15, 1).Value = "employeesListsPath"

# Row 19 name
$ws.Cells.Item(19, 1).Value = "emailMessageEmployeesList"

# Row 16 (formerly cvIndexPath) is removed entirely
$ws.Range("A16:C16").ClearContents()

# Row 17 now holds what used to live in row 18 (replyMessageNoAttachment)
$ws.Cells.Item(17, 1).Value = "replyMessageNoAttachment"
$ws.Cells.Item(17, 2).Value = "The mail you sent does not include an attachment."

# Row 18 now holds what used to live in row 19 (replyMessageConfirmation)
$ws.Cells.Item(18, 1).Value = "replyMessageConfirmation"
$ws.Cells.Item(18, 2).Value = "We have received the email sent by you."

# Row 21 (apiKey) remains unchanged.

# Update the active selection to match the saved view state
$ws.Range("B27").Select()
